$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LCOE")

# Update run_name values: "base" -> "base_case" (row 2), "base" -> "base_case_PV" (row 3)
$ws.Range("A2").Value = "base_case"
$ws.Range("A3").Value = "base_case_PV"

# Move the active selection from F16 to I16
$ws.Range("I16").Select()
